$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns for the crypto rows (2-51) with the
# latest scraped figures. D-column values whose new text parses as a plain number (e.g.
# "0.9999", "1.000") are written through a Text-formatted cell first, otherwise Excel
# would auto-convert the typed string into a float and silently drop significant
# trailing zeros; the cell format is then restored to the default "Normal" style so no
# extra formatting is left behind. Cells whose value does not actually change (e.g. D44,
# which keeps its old "1.000") are left untouched.

$ws.Range('D2').Value = '26.440.89'
$ws.Range('E2').Value = '  -0.44%  '

$ws.Range('D3').Value = '1.843.66'
$ws.Range('E3').Value = '  -1.60%  '

$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.9999'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.43%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '260.45'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -7.72%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.9999'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.28%  '

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5136'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '

$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3208'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -9.20%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.06758'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -1.60%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '18.78'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -6.65%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.7683'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -5.94%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.07683'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -1.00%  '

$ws.Range('D13').Value = '1.867.09'
$ws.Range('E13').Value = '  -0.38%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '88.76'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.70%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '5.025'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -1.87%  '

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +0.66%  '

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '14.06'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -1.84%  '

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.21%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.000007901'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -2.76%  '

$ws.Range('D20').Value = '26.466.29'
$ws.Range('E20').Value = '  -0.36%  '

$ws.Range('D21').Value = '2.100.39'
$ws.Range('E21').Value = '  +0.32%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '4.571'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -4.98%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '9.554'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -5.67%  '

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '5.955'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -4.58%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.324'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -3.10%  '

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '144.95'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +0.30%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '1.660'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -0.06%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '16.95'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -2.16%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '111.04'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.42%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '4.176'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -4.83%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '4.162'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -4.01%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.08731'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -0.74%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.04814'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -1.97%  '

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '1.134'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -3.56%  '

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.839'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.62%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.6864'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -7.89%  '

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '3.080'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -5.86%  '

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01807'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -3.34%  '

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.204'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -8.79%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.4911'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -6.08%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '113.33'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -2.77%  '

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.9042'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -6.59%  '

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '6.140'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -2.63%  '

$ws.Range('E44').Value = '  +0.41%  '

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '7.763'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -3.99%  '

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.4242'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -7.32%  '

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.1270'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -6.95%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '9.105'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -3.50%  '

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.05886'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.53%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '34.96'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -4.28%  '

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.421'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -6.10%  '

